# Remove form_id from basic forms: delete column B ("form_id") on the
# "settings" sheet, shifting version/style/namespaces left.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("settings")

# Capture the comment text that should survive (currently anchored at
# C1/D1/E1) before the column shift invalidates the layout.
$versionComment    = $ws.Range("C1").Comment.Text()
$styleComment      = $ws.Range("D1").Comment.Text()
$namespacesComment = $ws.Range("E1").Comment.Text()

# Delete the whole form_id column; cells to the right shift left
# automatically (shared strings / widths / dimension follow).
$ws.Columns.Item(2).Delete()

# Comments stay anchored to their original cell refs after a column
# delete in this automation surface, so re-home them by hand: overwrite
# the (now stale) comments in place so authorship/formatting carries
# over, then drop the now-duplicate trailing one.
$ws.Range("B1").Comment.Text($versionComment) | Out-Null
$ws.Range("C1").Comment.Text($styleComment) | Out-Null
$ws.Range("D1").Comment.Text($namespacesComment) | Out-Null
$ws.Range("E1").Comment.Delete() | Out-Null
